{"js": "// Replace each old division expression with its corresponding new one.\n// Every text run in this document's quiz table is unique, so a simple\n// search-and-replace keyed on the exact old text is unambiguous.\nconst replacements = [\n  [\"732\u00f76=\", \"679\u00f78=\"],\n  [\"925\u00f78=\", \"961\u00f78=\"],\n  [\"555\u00f79=\", \"676\u00f77=\"],\n  [\"634\u00f74=\", \"536\u00f76=\"],\n  [\"380\u00f78=\", \"616\u00f79=\"],\n  [\"852\u00f76=\", \"861\u00f78=\"],\n  [\"807\u00f74=\", \"469\u00f74=\"],\n  [\"178\u00f77=\", \"353\u00f78=\"],\n  [\"612\u00f78=\", \"366\u00f78=\"],\n  [\"998\u00f79=\", \"616\u00f76=\"],\n  [\"568\u00f78=\", \"827\u00f76=\"],\n  [\"428\u00f73=\", \"466\u00f74=\"],\n  [\"643\u00f75=\", \"369\u00f79=\"],\n  [\"600\u00f76=\", \"981\u00f75=\"],\n  [\"270\u00f73=\", \"484\u00f74=\"],\n  [\"848\u00f79=\", \"801\u00f72=\"],\n  [\"449\u00f72=\", \"689\u00f78=\"],\n  [\"701\u00f76=\", \"679\u00f73=\"],\n  [\"323\u00f74=\", \"246\u00f75=\"],\n  [\"334\u00f72=\", \"359\u00f73=\"],\n  [\"482\u00f75=\", \"326\u00f74=\"],\n  [\"676\u00f72=\", \"325\u00f76=\"],\n  [\"114\u00f79=\", \"183\u00f75=\"],\n  [\"506\u00f75=\", \"961\u00f75=\"],\n  [\"328\u00f73=\", \"637\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old division expression with its corresponding new one.\n# Every source string in this quiz table is unique, so Find/Replace keyed\n# on the exact old text (with wdReplaceAll) is unambiguous and idempotent.\n$replacements = @(\n    @(\"732\u00f76=\", \"679\u00f78=\"),\n    @(\"925\u00f78=\", \"961\u00f78=\"),\n    @(\"555\u00f79=\", \"676\u00f77=\"),\n    @(\"634\u00f74=\", \"536\u00f76=\"),\n    @(\"380\u00f78=\", \"616\u00f79=\"),\n    @(\"852\u00f76=\", \"861\u00f78=\"),\n    @(\"807\u00f74=\", \"469\u00f74=\"),\n    @(\"178\u00f77=\", \"353\u00f78=\"),\n    @(\"612\u00f78=\", \"366\u00f78=\"),\n    @(\"998\u00f79=\", \"616\u00f76=\"),\n    @(\"568\u00f78=\", \"827\u00f76=\"),\n    @(\"428\u00f73=\", \"466\u00f74=\"),\n    @(\"643\u00f75=\", \"369\u00f79=\"),\n    @(\"600\u00f76=\", \"981\u00f75=\"),\n    @(\"270\u00f73=\", \"484\u00f74=\"),\n    @(\"848\u00f79=\", \"801\u00f72=\"),\n    @(\"449\u00f72=\", \"689\u00f78=\"),\n    @(\"701\u00f76=\", \"679\u00f73=\"),\n    @(\"323\u00f74=\", \"246\u00f75=\"),\n    @(\"334\u00f72=\", \"359\u00f73=\"),\n    @(\"482\u00f75=\", \"326\u00f74=\"),\n    @(\"676\u00f72=\", \"325\u00f76=\"),\n    @(\"114\u00f79=\", \"183\u00f75=\"),\n    @(\"506\u00f75=\", \"961\u00f75=\"),\n    @(\"328\u00f73=\", \"637\u00f78=\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
